# Update of GINF1 Modules
# Fill in the GINF1 module table (rows 2-13) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "GINF11"
$ws.Cells.Item(2, 2).Value = "Maths pour l'ingénieur"
$ws.Cells.Item(2, 4).Value = "Statistiques"

$ws.Cells.Item(3, 1).Value = "GINF12"
$ws.Cells.Item(3, 2).Value = "Signal"
$ws.Cells.Item(3, 3).Value = "Siham"
$ws.Cells.Item(3, 4).Value = "Traitement de Signal"

$ws.Cells.Item(4, 1).Value = "GINF13"
$ws.Cells.Item(4, 2).Value = "Electronique 1"
$ws.Cells.Item(4, 3).Value = "El Oualkadi"
$ws.Cells.Item(4, 4).Value = "Electronique analogique, Electronique numérique"

$ws.Cells.Item(5, 1).Value = "GINF14"
$ws.Cells.Item(5, 2).Value = "Programmation"
$ws.Cells.Item(5, 3).Value = "Amechnoue"
$ws.Cells.Item(5, 4).Value = "C, Programmation web PHP"

$ws.Cells.Item(6, 1).Value = "GINF15"
$ws.Cells.Item(6, 2).Value = "BD & Réseaux"
$ws.Cells.Item(6, 3).Value = "Tanana"
$ws.Cells.Item(6, 4).Value = "BD relationnelle, Concept fondamentaux des réseaux"

$ws.Cells.Item(7, 1).Value = "GINF16"
$ws.Cells.Item(7, 2).Value = "Langues et Communication"
$ws.Cells.Item(7, 3).Value = "Haris"
$ws.Cells.Item(7, 4).Value = "Développement personnel"

$ws.Cells.Item(8, 1).Value = "GINF21"
$ws.Cells.Item(8, 2).Value = "Développement Informatique"
$ws.Cells.Item(8, 3).Value = "El Haddad"
$ws.Cells.Item(8, 4).Value = "POO C++,Programmation Web PHP5 & Mysql"

$ws.Cells.Item(9, 1).Value = "GINF22"
$ws.Cells.Item(9, 2).Value = "BD : Développemenent et Manipulation"
$ws.Cells.Item(9, 3).Value = "Fissoune Rachida"
$ws.Cells.Item(9, 4).Value = "Méthodes et modélisation BD, PLSQL"

$ws.Cells.Item(10, 1).Value = "GINF23"
$ws.Cells.Item(10, 2).Value = "Théories des Graphes et recherches opérationnelle"
$ws.Cells.Item(10, 3).Value = "Samadi"
$ws.Cells.Item(10, 4).Value = "Recherche Opérationnel, Théorie des graphes"

$ws.Cells.Item(11, 1).Value = "GINF24"
$ws.Cells.Item(11, 2).Value = "Réseaux 1"
$ws.Cells.Item(11, 3).Value = "Tanana"
$ws.Cells.Item(11, 4).Value = "Protocole et Adressage réseaux, Technologie des réseaux, TP CISCO"

$ws.Cells.Item(12, 1).Value = "GINF25"
$ws.Cells.Item(12, 2).Value = "Architectures et Linux"
$ws.Cells.Item(12, 3).Value = "Belmokadem"
$ws.Cells.Item(12, 4).Value = "Micro-Achitecture des processeur, Assembleur, Linux"

$ws.Cells.Item(13, 1).Value = "GINF26"
$ws.Cells.Item(13, 2).Value = "Management de l'entreprise 1"
$ws.Cells.Item(13, 3).Value = "NAIT BOUKER Nezha"
$ws.Cells.Item(13, 4).Value = "Gestion des entreprises, Comptabilité, Economie"

# Column widths (B, C, D) to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 50.33333333333333
$ws.Columns.Item(3).ColumnWidth = 25.0
$ws.Columns.Item(4).ColumnWidth = 47.33333333333333

# Selection moves to C8 (as in the committed workbook).
[void]$ws.Range("C8").Select()
